# 9th Stab - Cosmetic Changes
# Shift the existing 10 days of data (columns B:K) three columns to the
# right (columns E:N) to make room for three newly-scraped days, then
# populate the freed-up B:D columns with the new day's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Shift existing data (B1:K6) three columns right, into E1:N6 ---
$src = $ws.Range("B1:K6")
$dst = $ws.Range("E1:N6")
$dst.Value2 = $src.Value2

# --- 2. Populate the new B:D columns with the latest scraped data -----
# Row 1: dates
$ws.Range("B1").Value2 = "Jun_18"
$ws.Range("C1").Value2 = "Jun_17"
$ws.Range("D1").Value2 = "Jun_15"

# Row 2: IRT
$ws.Range("B2").Value2 = "Hold       (`$10.17)"
$ws.Range("C2").Value2 = "Hold       (`$10.15)"
$ws.Range("D2").Value2 = "Hold       (`$10.15)"

# Row 3: HCP
$ws.Range("B3").Value2 = "Hold       (`$23.95)"
$ws.Range("C3").Value2 = "Hold       (`$24.02)"
$ws.Range("D3").Value2 = "Hold       (`$24.02)"

# Row 4: KIM
$ws.Range("B4").Value2 = "Hold       (`$16.36)"
$ws.Range("C4").Value2 = "Hold       (`$16.50)"
$ws.Range("D4").Value2 = "Hold       (`$16.50)"

# Row 5: MORT (unchanged value, written explicitly for completeness)
$ws.Range("B5").Value2 = "UN         (0)"
$ws.Range("C5").Value2 = "UN         (0)"
$ws.Range("D5").Value2 = "UN         (0)"

# Row 6: PLD
$ws.Range("B6").Value2 = "Buy        (`$64.17)"
$ws.Range("C6").Value2 = "Buy        (`$63.96)"
$ws.Range("D6").Value2 = "Buy        (`$63.96)"
